# Refresh the cryptocurrency price/volume snapshot (GitHub Actions cron job).
#
# Column D ("Price") holds text like "1.00" or "592.02" that *looks* numeric.
# Assigning such a string straight to .Value lets Excel auto-convert it to a
# real number (dropping significant trailing zeros / the original text
# formatting), so for those cells we briefly force the cell to Text format,
# write the value, then restore the "Normal" style so no stray number-format
# change is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "70.263.53"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +5.29%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "3.612.17"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +5.15%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.03%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "592.02"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +3.94%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "189.64"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +4.10%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.646"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +2.28%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "3.605.66"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +5.16%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.01%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.178"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +3.08%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.662"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +3.01%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "58.35"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +4.34%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000288"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +3.47%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "9.89"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +5.64%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "4.189.37"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +5.13%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "19.59"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +5.61%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "3.613.36"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +5.18%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "70.227.91"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +5.12%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "12.53"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +4.15%  "
$ws.Cells.Item(20, 5).Value = "  +0.54%  "
$ws.Cells.Item(21, 5).Value = "  +3.99%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "491.22"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.27%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "17.79"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +10.03%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "5.37"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +7.47%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "4.47"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +6.54%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "90.77"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.81%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "3.12"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +5.16%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "11.18"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +2.06%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "9.36"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +3.81%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "32.70"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +4.33%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "7.62"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +5.18%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "12.33"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +5.85%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "619.92"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +4.54%  "
$ws.Cells.Item(34, 5).Value = "  +6.62%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "65.56"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +4.02%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0826"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +7.20%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "38.38"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +5.91%  "
$ws.Cells.Item(38, 2).Value = "Kaspa"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.147"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.94%  "
$ws.Cells.Item(39, 5).Value = "  +0.12%  "
$ws.Cells.Item(40, 2).Value = "TheGraph"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.402"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +3.63%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "3.58"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.05%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "3.322.58"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +5.69%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "3.16"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +8.60%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.0454"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +6.19%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "2.72"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +6.54%  "
$ws.Cells.Item(46, 2).Value = "Stellar"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +2.09%  "
$ws.Cells.Item(47, 2).Value = "ApeXProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "3.27"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.15%  "
$ws.Cells.Item(48, 2).Value = "dogwifhat"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "2.73"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -2.42%  "
$ws.Cells.Item(49, 2).Value = "THORChain"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "9.06"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.65%  "
$ws.Cells.Item(50, 5).Value = "  +5.54%  "
$ws.Cells.Item(51, 2).Value = "Monero"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "142.86"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.47%  "
